$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New default-item rows (6-17) -- shared strings are interned in the order
# they're first referenced, matching the target sharedStrings.xml ordering.

# SALUNDRA
$ws.Range("A6").Value = "SALUNDRA"
$ws.Range("B6").Value = "Fabulous Hat"
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = "SALUNDRA"
$ws.Range("B7").Value = "Clothes"
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = "SALUNDRA"
$ws.Range("B8").Value = "Uniform"
$ws.Range("C8").Value = 1

# MOLRELLA
$ws.Range("A9").Value = "MOLRELLA"
$ws.Range("B9").Value = "Healing Draught"
$ws.Range("C9").Value = 0

$ws.Range("A10").Value = "MOLRELLA"
$ws.Range("B10").Value = "Lock Picks"
$ws.Range("C10").Value = 0

# FERDINAND
$ws.Range("A11").Value = "FERDINAND"
$ws.Range("B11").Value = "Grimoire"
$ws.Range("C11").Value = 0

$ws.Range("A12").Value = "FERDINAND"
$ws.Range("B12").Value = "Amethyst Wizard's Uniform"
$ws.Range("C12").Value = 2

$ws.Range("A13").Value = "FERDINAND"
$ws.Range("B13").Value = "6 sheets of Parchment"
$ws.Range("C13").Value = 0

$ws.Range("A14").Value = "FERDINAND"
$ws.Range("B14").Value = "Quill and Ink"
$ws.Range("C14").Value = 0

# AMRIS
$ws.Range("A15").Value = "AMRIS"
$ws.Range("B15").Value = "Healing Draught"
$ws.Range("C15").Value = 0

$ws.Range("A16").Value = "AMRIS"
$ws.Range("B16").Value = "High Elf Clothing"
$ws.Range("C16").Value = 1

# ELSE
$ws.Range("A17").Value = "ELSE"
$ws.Range("B17").Value = "10 Bullets"
$ws.Range("C17").Value = 2

# Column widths widened for the longer character/item names now present.
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 23.666666666666668

# Selection moves to below the newly added data, matching the author's
# last-used cell after populating the sheet.
$ws.Range("B19").Select()
